$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price/Volume(1h) columns in this sheet always hold plain text (e.g.
# the thousands-dotted "97.389.78" or the padded "  +3.14%  "), never real
# numbers. Pin every cell we're about to rewrite to the Text number format
# first (cell by cell - a Union range only honours the first area) so Excel's
# automatic type inference can't quietly turn a value like "24.70" into the
# shorter-looking number 24.7.
$textCells = @("D2","E2","D3","E3","E4","D5","E5","D6","E6","E7","E8","E9","D10","E10","D11","E11","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","E29","E30","D31","E31","E32","E33","D34","E34","D35","E35","D36","E36","D37","E37","D38","E38","D39","E39","D40","E40","D41","E41","D42","E42","D43","E43","D44","E44","E45","E46","D47","E47","D48","E48","D49","E49","D50","E50","D51","E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '97.389.78'
$ws.Range("E2").Value = '  +3.14%  '

# Row 3
$ws.Range("D3").Value = '3.362.71'
$ws.Range("E3").Value = '  +6.76%  '

# Row 4
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$ws.Range("D5").Value = '246.53'
$ws.Range("E5").Value = '  +2.26%  '

# Row 6
$ws.Range("D6").Value = '626.87'
$ws.Range("E6").Value = '  +1.25%  '

# Row 7
$ws.Range("E7").Value = '  -0.24%  '

# Row 8
$ws.Range("E8").Value = '  -1.17%  '

# Row 9
$ws.Range("E9").Value = '  +0.05%  '

# Row 10
$ws.Range("D10").Value = '3.362.81'
$ws.Range("E10").Value = '  +7.11%  '

# Row 11
$ws.Range("D11").Value = '0.798'
$ws.Range("E11").Value = '  -4.68%  '

# Row 12
$ws.Range("E12").Value = '  +0.78%  '

# Row 13
$ws.Range("D13").Value = '97.381.84'
$ws.Range("E13").Value = '  +3.52%  '

# Row 14
$ws.Range("D14").Value = '0.0000252'
$ws.Range("E14").Value = '  +2.12%  '

# Row 15
$ws.Range("D15").Value = '35.92'
$ws.Range("E15").Value = '  +1.61%  '

# Row 16
$ws.Range("D16").Value = '3.963.94'
$ws.Range("E16").Value = '  +6.59%  '

# Row 17
$ws.Range("D17").Value = '5.54'
$ws.Range("E17").Value = '  +2.08%  '

# Row 18
$ws.Range("D18").Value = '3.345.00'
$ws.Range("E18").Value = '  +8.22%  '

# Row 19
$ws.Range("E19").Value = '  -2.43%  '

# Row 20
$ws.Range("D20").Value = '15.35'
$ws.Range("E20").Value = '  +0.92%  '

# Row 21
$ws.Range("D21").Value = '491.05'
$ws.Range("E21").Value = '  +8.66%  '

# Row 22
$ws.Range("D22").Value = '0.0000214'
$ws.Range("E22").Value = '  +4.81%  '

# Row 23
$ws.Range("D23").Value = '5.94'
$ws.Range("E23").Value = '  -1.09%  '

# Row 24
$ws.Range("D24").Value = '9.43'
$ws.Range("E24").Value = '  +2.75%  '

# Row 25
$ws.Range("D25").Value = '5.76'
$ws.Range("E25").Value = '  +0.24%  '

# Row 26
$ws.Range("D26").Value = '88.65'
$ws.Range("E26").Value = '  +2.57%  '

# Row 27
$ws.Range("D27").Value = '12.21'
$ws.Range("E27").Value = '  -0.78%  '

# Row 28
$ws.Range("D28").Value = '3.511.26'
$ws.Range("E28").Value = '  +6.88%  '

# Row 29
$ws.Range("E29").Value = '  +0.07%  '

# Row 30
$ws.Range("E30").Value = '  +1.23%  '

# Row 31
$ws.Range("D31").Value = '0.242'
$ws.Range("E31").Value = '  -6.10%  '

# Row 32
$ws.Range("E32").Value = '  -2.83%  '

# Row 33
$ws.Range("E33").Value = '  +0.03%  '

# Row 34
$ws.Range("D34").Value = '9.45'
$ws.Range("E34").Value = '  -0.05%  '

# Row 35
$ws.Range("D35").Value = '27.87'
$ws.Range("E35").Value = '  +5.11%  '

# Row 36
$ws.Range("D36").Value = '0.156'
$ws.Range("E36").Value = '  -4.49%  '

# Row 37
$ws.Range("D37").Value = '7.56'
$ws.Range("E37").Value = '  -6.53%  '

# Row 38
$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '1.97'
$ws.Range("E38").Value = '  +2.14%  '

# Row 39
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D39").Value = '507.16'
$ws.Range("E39").Value = '  +4.22%  '

# Row 40
$ws.Range("B40").Value = 'PolygonEcosystemToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D40").Value = '0.459'
$ws.Range("E40").Value = '  -0.49%  '

# Row 41
$ws.Range("B41").Value = 'WhiteBITCoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D41").Value = '24.70'
$ws.Range("E41").Value = '  +3.11%  '

# Row 42
$ws.Range("D42").Value = '1.30'
$ws.Range("E42").Value = '  -0.32%  '

# Row 43
$ws.Range("D43").Value = '0.816'
$ws.Range("E43").Value = '  +15.13%  '

# Row 44
$ws.Range("D44").Value = '3.33'
$ws.Range("E44").Value = '  +0.05%  '

# Row 45
$ws.Range("E45").Value = '  -7.35%  '

# Row 46
$ws.Range("E46").Value = '  +0.03%  '

# Row 47
$ws.Range("D47").Value = '161.33'
$ws.Range("E47").Value = '  +1.58%  '

# Row 48
$ws.Range("D48").Value = '1.97'
$ws.Range("E48").Value = '  +4.23%  '

# Row 49
$ws.Range("D49").Value = '4.62'
$ws.Range("E49").Value = '  +2.66%  '

# Row 50
$ws.Range("B50").Value = 'Optimism'
$ws.Range("C50").Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range("D50").Value = '2.19'
$ws.Range("E50").Value = '  +21.95%  '

# Row 51
$ws.Range("D51").Value = '1.37'
$ws.Range("E51").Value = '  +1.79%  '
